$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (values look like plain decimals).
$textCells = @("D5","D6","D10","D12","D16","D17","D19","D20","D24","D29","D32","D33","D35","D36","D37","D38","D40","D44","D45","D47","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the latest cryptos data pull
$ws.Range("D2").Value = "55.322.70"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "2.353.91"
$ws.Range("E3").Value = "  -4.39%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "477.36"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").Value = "146.49"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +20.64%  "
$ws.Range("D9").Value = "2.359.09"
$ws.Range("E9").Value = "  -4.44%  "
$ws.Range("D10").Value = "0.0963"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  -5.61%  "
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "2.764.91"
$ws.Range("E14").Value = "  -4.51%  "
$ws.Range("D15").Value = "55.235.32"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "19.99"
$ws.Range("E16").Value = "  -4.57%  "
$ws.Range("D17").Value = "0.0000130"
$ws.Range("E17").Value = "  -4.20%  "
$ws.Range("D18").Value = "2.354.28"
$ws.Range("E18").Value = "  -4.87%  "
$ws.Range("D19").Value = "4.59"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").Value = "314.51"
$ws.Range("E21").Value = "  -4.44%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").Value = "56.44"
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -3.35%  "
$ws.Range("E27").Value = "  -4.43%  "
$ws.Range("D28").Value = "2.451.61"
$ws.Range("E28").Value = "  -4.76%  "
$ws.Range("D29").Value = "7.07"
$ws.Range("E29").Value = "  -7.84%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "0.0₃0745"
$ws.Range("E31").Value = "  -4.57%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "18.16"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "145.59"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").Value = "5.11"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "3.61"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.10"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("D38").Value = "0.811"
$ws.Range("E38").Value = "  -5.10%  "
$ws.Range("E39").Value = "  +10.59%  "
$ws.Range("D40").Value = "33.71"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("D44").Value = "0.578"
$ws.Range("E44").Value = "  -3.54%  "
$ws.Range("D45").Value = "0.0517"
$ws.Range("E45").Value = "  -5.92%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "249.68"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("E49").Value = "  -6.51%  "
$ws.Range("D50").Value = "1.798.94"
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("D51").Value = "16.63"
$ws.Range("E51").Value = "  -4.85%  "
